# chore: adapt column header formatting to respective input file names
#
# Renames the "_old" / "_new" header-name suffixes to the concrete
# format-version identifiers ("_FV2310" / "_FV2404"), turns row 1 into a
# frozen header row, and wraps the data range in a native Excel table
# (ListObject) with an AutoFilter, matching the new AHB-Diff export shape.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename the header cells (A1:J1 = "FV2310" / "before" columns, -----
#        L1:U1 = "FV2404" / "after" columns; K1 "diff" is unchanged) -------
$oldHeaders = @(
    "Segmentname",
    "Segmentgruppe",
    "Segment",
    "Datenelement",
    "Segment ID",
    "Code",
    "Qualifier",
    "Beschreibung",
    "Bedingungsausdruck",
    "Bedingung"
)

$col = 1
foreach ($name in $oldHeaders) {
    $ws.Cells.Item(1, $col).Value = "$($name)_FV2310"
    $col++
}

# Column K (11) holds the "diff" header and stays untouched.

$newHeaders = $oldHeaders
$col = 12
foreach ($name in $newHeaders) {
    $ws.Cells.Item(1, $col).Value = "$($name)_FV2404"
    $col++
}

# --- 2. Freeze the header row (split after row 1) --------------------------
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

# --- 3. Turn the used range into a table with an autofilter ----------------
$rng = $ws.Range("A1:U80")
$tbl = $ws.ListObjects.Add(1, $rng, [System.Type]::Missing, 1)
$tbl.Name = "Table1"
$tbl.TableStyle = ""
